$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "L5FzvRYyZyRtrIofAyYS"
$ws.Range("A3").Value = "wMVY6Ls9uzob9kch569W"
